$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "upper95cl"
$ws.Range("H1").Value = "lower998cl"
$ws.Range("F1").Value = "lower95cl"
$ws.Range("I1").Value = "upper998cl"

$ws.Range("F8").Select()
